$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the existing Synthetic-0.5 block (rows 32-37)
# down into the two new blocks (rows 38-43 and 44-49) so that column A retains
# its bold/border/center style, matching the pattern already used in the sheet.
$ws.Range("A32:I37").Copy($ws.Range("A38"))
$ws.Range("A32:I37").Copy($ws.Range("A44"))

# Now overwrite the copied values with the actual new experiment results.
$ws.Range("A38").Value = 0
$ws.Range("B38").Value = "Synthetic-0.5"
$ws.Range("C38").Value = "LGBM Baseline"
$ws.Range("D38").Value = 0.565989847715736
$ws.Range("E38").Value = 0.565989847715736
$ws.Range("F38").Value = 0.9968309859154929
$ws.Range("G38").Value = 0.9968309859154929
$ws.Range("H38").Value = 0.9981804369560553
$ws.Range("I38").Value = 0.9981804369560553

$ws.Range("A39").Value = 1
$ws.Range("B39").Value = "Synthetic-0.5"
$ws.Range("C39").Value = "LGBM Upsample"
$ws.Range("D39").Value = 0.1715260545905707
$ws.Range("E39").Value = 0.9748037512393462
$ws.Range("F39").Value = 0.9505003706449222
$ws.Range("G39").Value = 0.9742003986364398
$ws.Range("H39").Value = 0.994321693648723
$ws.Range("I39").Value = 0.9943616762474357

$ws.Range("A40").Value = 2
$ws.Range("B40").Value = "Synthetic-0.5"
$ws.Range("C40").Value = "LGBM Downsample"
$ws.Range("D40").Value = 0.03798162621692033
$ws.Range("E40").Value = 1
$ws.Range("F40").Value = 0.7399555226093403
$ws.Range("G40").Value = 1
$ws.Range("H40").Value = 0.9711692506168886
$ws.Range("I40").Value = 1

$ws.Range("A41").Value = 3
$ws.Range("B41").Value = "Synthetic-0.5"
$ws.Range("C41").Value = "SMOTE LGBM"
$ws.Range("D41").Value = 0.07346881765486417
$ws.Range("E41").Value = 0.9143096339772789
$ws.Range("F41").Value = 0.8770663454410674
$ws.Range("G41").Value = 0.9113406478773541
$ws.Range("H41").Value = 0.9708859657552483
$ws.Range("I41").Value = 0.9692645787690264

$ws.Range("A42").Value = 4
$ws.Range("B42").Value = "Synthetic-0.5"
$ws.Range("C42").Value = "LGBM Balanced Bagging"
$ws.Range("D42").Value = 0.3777452415812592
$ws.Range("E42").Value = 0.3777452415812592
$ws.Range("F42").Value = 0.9960618977020015
$ws.Range("G42").Value = 0.9960618977020015
$ws.Range("H42").Value = 0.9999730501107887
$ws.Range("I42").Value = 0.9999730501107887

$ws.Range("A43").Value = 5
$ws.Range("B43").Value = "Synthetic-0.5"
$ws.Range("C43").Value = "LGBM_Imbalance"
$ws.Range("D43").Value = 0.0954427433073686
$ws.Range("E43").Value = 0.0954427433073686
$ws.Range("F43").Value = 0.9063843587842847
$ws.Range("G43").Value = 0.9063843587842847
$ws.Range("H43").Value = 0.9576621527591677
$ws.Range("I43").Value = 0.9576621527591677

$ws.Range("A44").Value = 0
$ws.Range("B44").Value = "Synthetic-0.5"
$ws.Range("C44").Value = "LGBM Baseline"
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0.5572916666666667
$ws.Range("F44").Value = 0.9947741002927987
$ws.Range("G44").Value = 0.9968495181616012
$ws.Range("H44").Value = 0.7856255670059564
$ws.Range("I44").Value = 0.9989981110762746

$ws.Range("A45").Value = 1
$ws.Range("B45").Value = "Synthetic-0.5"
$ws.Range("C45").Value = "LGBM Upsample"
$ws.Range("D45").Value = 0.0648854961832061
$ws.Range("E45").Value = 0.9742856752684549
$ws.Range("F45").Value = 0.945517215818539
$ws.Range("G45").Value = 0.9736942592109675
$ws.Range("H45").Value = 0.8444959591937506
$ws.Range("I45").Value = 0.9939958935087631

$ws.Range("A46").Value = 2
$ws.Range("B46").Value = "Synthetic-0.5"
$ws.Range("C46").Value = "LGBM Downsample"
$ws.Range("D46").Value = 0.02641408751334045
$ws.Range("E46").Value = 1
$ws.Range("F46").Value = 0.729513361254216
$ws.Range("G46").Value = 1
$ws.Range("H46").Value = 0.8133324028620204
$ws.Range("I46").Value = 1

$ws.Range("A47").Value = 3
$ws.Range("B47").Value = "Synthetic-0.5"
$ws.Range("C47").Value = "SMOTE LGBM"
$ws.Range("D47").Value = 0.04306350504224585
$ws.Range("E47").Value = 0.9065770245264153
$ws.Range("F47").Value = 0.8698713909788369
$ws.Range("G47").Value = 0.9038436463882576
$ws.Range("H47").Value = 0.8431981862846573
$ws.Range("I47").Value = 0.9655339671578294

$ws.Range("A48").Value = 4
$ws.Range("B48").Value = "Synthetic-0.5"
$ws.Range("C48").Value = "LGBM Balanced Bagging"
$ws.Range("D48").Value = 0
$ws.Range("E48").Value = 0.4069767441860465
$ws.Range("F48").Value = 0.9952929839516697
$ws.Range("G48").Value = 0.9962194217939214
$ws.Range("H48").Value = 0.8249686698971221
$ws.Range("I48").Value = 0.9999287727692603

$ws.Range("A49").Value = 5
$ws.Range("B49").Value = "Synthetic-0.5"
$ws.Range("C49").Value = "LGBM_Imbalance"
$ws.Range("D49").Value = 0.02310924369747899
$ws.Range("E49").Value = 0.05368421052631579
$ws.Range("F49").Value = 0.8621251992142619
$ws.Range("G49").Value = 0.8667160859896219
$ws.Range("H49").Value = 0.6034083985200814
$ws.Range("I49").Value = 0.8793823951574593
